$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting the existing rows 48-49 down to 49-50.
$ws.Rows.Item(48).Insert()

# Fill the newly inserted row 48 with the new weekly price entry.
$ws.Cells.Item(48, 1).Value = 1
$ws.Cells.Item(48, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(48, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(48, 4).Value = 44931
$ws.Cells.Item(48, 5).Value = 15
$ws.Cells.Item(48, 6).Value = 100112052
$ws.Cells.Item(48, 7).Value = "Albahaca"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Segunda"
$ws.Cells.Item(48, 10).Value = 300
$ws.Cells.Item(48, 11).Value = 1300
$ws.Cells.Item(48, 12).Value = 1500
$ws.Cells.Item(48, 13).Value = 1400
$ws.Cells.Item(48, 14).Value = "$/paquete"
$ws.Cells.Item(48, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(48, 16).Value = 1400
$ws.Cells.Item(48, 17).Value = 1
$ws.Cells.Item(48, 18).Value = "Hortaliza"
